$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10 (columns G-T)
$values = @{
  2 = @{ G=0.3997546666666666; H=1.199264; I=0.1320462084214824; J=0.1320462084214824; M=2.341355666666667; N=7.024067000000001; O=0.03973512964576821; P=0.0397351296457682; Q=0.9359678540764444; R=8.423710686688; S=0.005246873210859733; T=0.005246873210859732 }
  3 = @{ G=0.3997546666666666; H=1.199264; I=0.1320462084214824; J=0.1320462084214824; O=0.5779093692199981; P=0.5779093692199981; Q=13.61275518619378; R=122.514796675744; S=0.0763107410167513; T=0.0763107410167513 }
  4 = @{ G=0.3997546666666666; H=1.199264; I=0.1320462084214824; J=0.1320462084214824; O=0.3823555011342337; P=0.3823555011342337; Q=9.006449987235555; R=81.05804988512; S=0.05048859419387137; T=0.05048859419387137 }
  5 = @{ I=0.6840925621829359; J=0.684092562182936; M=2.341355666666667; N=7.024067000000001; O=0.03973512964576821; P=0.0397351296457682; Q=4.848974121030889; R=43.64076708927801; S=0.02718250664804471; T=0.02718250664804471 }
  6 = @{ I=0.6840925621829359; J=0.684092562182936; O=0.5779093692199981; P=0.5779093692199981; Q=70.52368019510156; R=634.713121755914; S=0.3953435010992328; T=0.3953435010992329 }
  7 = @{ I=0.6840925621829359; J=0.684092562182936; O=0.3823555011342337; P=0.3823555011342337; S=0.2615665544356583; T=0.2615665544356584 }
  8 = @{ G=0.5566186666666667; I=0.1838612293955817; J=0.1838612293955817; M=2.341355666666667; N=7.024067000000001; O=0.03973512964576821; P=0.0397351296457682; Q=1.303242269372445; R=11.729180424352; S=0.00730574978686377; T=0.007305749786863768 }
  9 = @{ G=0.5566186666666667; I=0.1838612293955817; J=0.1838612293955817; O=0.5779093692199981; P=0.5779093692199981; S=0.106255127104014; T=0.106255127104014 }
  10 = @{ G=0.5566186666666667; I=0.1838612293955817; J=0.1838612293955817; O=0.3823555011342337; P=0.3823555011342337; S=0.07030035250470396; T=0.07030035250470394 }
}

foreach ($row in $values.Keys) {
  $cols = $values[$row]
  foreach ($col in $cols.Keys) {
    $addr = "$col$row"
    $ws.Range($addr).Value = $cols[$col]
  }
}
